$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 gets the values that used to be row 16's scraped data
$ws.Range("A15").Value = "IT Support Technician 1"
$ws.Range("C15").Value = "Not Disclosed"
$ws.Range("E15").Value = "Full-time"

# Row 16 gets the newly scraped job data
$ws.Range("A16").Value = "IT Support Specialist Distribution Center"
$ws.Range("E16").Value = "$20.05 - $30.10 an hour"
